# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# These two sheets contain duplicate data sets, both need the same updates.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 317
    3  = 13923
    5  = 101
    7  = 288
    8  = 504
    9  = 12
    15 = 5935
    17 = 97
    18 = 988
    19 = 136
    22 = 303
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
